# Applies the cryptos-list price/volume refresh described in the commit diff.
# Cells in columns D/E hold plain display text (inlineStr in the source file),
# even when the text looks like a number (e.g. "601.84"). Assigning such a
# string straight to .Value would let Excel auto-convert the cell to a real
# number, which would NOT match the original text-formatted cells. For any
# new value that parses as a plain number we therefore force the cell to
# Text format first, write the value, then reset the style back to Normal so
# no stray number-format / quote-prefix style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    if ($Text -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $Cell.NumberFormat = "@"
        $Cell.Value = $Text
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Text
    }
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.908.34"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.515.67"
$ws.Range("E3").Value = "  -0.71%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "601.84"
$ws.Range("E5").Value = "  -1.61%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "196.09"
$ws.Range("E6").Value = "  +6.01%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.77%  "
# Row 8 - USDC
$ws.Range("E8").Value = "  -0.07%  "
# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.210"
$ws.Range("E9").Value = "  -2.26%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  +1.38%  "
# Row 11 - Avalanche
Set-TextValue $ws.Range("D11") "54.06"
$ws.Range("E11").Value = "  +1.12%  "

# Row 12 - ShibaInu
Set-TextValue $ws.Range("D12") "0.0000301"
$ws.Range("E12").Value = "  -2.68%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +1.04%  "
# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.068.91"
$ws.Range("E14").Value = "  -0.80%  "

# Row 15 - BitcoinCash
Set-TextValue $ws.Range("D15") "603.65"
$ws.Range("E15").Value = "  -0.84%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "70.025.29"
$ws.Range("E16").Value = "  +0.17%  "

# Row 17 - Chainlink
Set-TextValue $ws.Range("D17") "19.07"
$ws.Range("E17").Value = "  +1.24%  "

# Row 18 - Uniswap
$ws.Range("E18").Value = "  -0.55%  "
# Row 19 - WrappedEther
$ws.Range("D19").Value = "3.520.17"
$ws.Range("E19").Value = "  -1.48%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +0.76%  "
# Row 21 - Polygon
Set-TextValue $ws.Range("D21") "0.992"
$ws.Range("E21").Value = "  +0.14%  "

# Row 22 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D22") "18.21"
$ws.Range("E22").Value = "  +3.93%  "

# Row 23/24 - Toncoin and Litecoin swap ranking order
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D23") "5.23"
$ws.Range("E23").Value = "  +6.23%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D24") "103.74"
$ws.Range("E24").Value = "  +3.82%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "4.59"
$ws.Range("E25").Value = "  -2.76%  "

# Row 26 - ImmutableX
Set-TextValue $ws.Range("D26") "3.07"
$ws.Range("E26").Value = "  +2.78%  "

# Row 27 - RenderToken
Set-TextValue $ws.Range("D27") "10.92"
$ws.Range("E27").Value = "  -0.30%  "

# Row 28 - Filecoin
Set-TextValue $ws.Range("D28") "9.66"
$ws.Range("E28").Value = "  +0.94%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "33.53"
$ws.Range("E29").Value = "  +3.43%  "

# Row 30 - dogwifhat
Set-TextValue $ws.Range("D30") "4.53"
$ws.Range("E30").Value = "  +25.26%  "

# Row 31 - NEARProtocol
Set-TextValue $ws.Range("D31") "7.12"
$ws.Range("E31").Value = "  +1.48%  "

# Row 32 - Cosmos
Set-TextValue $ws.Range("D32") "12.72"
$ws.Range("E32").Value = "  +4.28%  "

# Row 33 - Hedera
Set-TextValue $ws.Range("D33") "0.116"
$ws.Range("E33").Value = "  +1.69%  "

# Row 34 - OKB
Set-TextValue $ws.Range("D34") "63.17"
$ws.Range("E34").Value = "  -0.41%  "

# Row 35 - Maker
$ws.Range("D35").Value = "3.742.73"
$ws.Range("E35").Value = "  +5.74%  "

# Row 36 - PEPE
$ws.Range("D36").Value = "0.0₃0817"
$ws.Range("E36").Value = "  +4.90%  "

# Row 37 - Dai
$ws.Range("E37").Value = "  +0.06%  "
# Row 39 - TheGraph
Set-TextValue $ws.Range("D39") "0.393"
$ws.Range("E39").Value = "  -1.94%  "

# Row 40/41 - Stacks and InjectiveProtocol swap ranking order
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D40") "3.59"
$ws.Range("E40").Value = "  +0.81%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D41") "36.78"
$ws.Range("E41").Value = "  -0.63%  "

# Row 42 - Bittensor
Set-TextValue $ws.Range("D42") "490.99"
$ws.Range("E42").Value = "  -8.01%  "

# Row 43 - Kaspa
$ws.Range("E43").Value = "  -0.19%  "
# Row 44 - VeChain
Set-TextValue $ws.Range("D44") "0.0456"
$ws.Range("E44").Value = "  +0.11%  "

# Row 45 - ApeXProtocol
$ws.Range("E45").Value = "  -0.72%  "
# Row 46 - Stellar
$ws.Range("E46").Value = "  -1.22%  "
# Row 47 - ThetaToken
Set-TextValue $ws.Range("D47") "2.83"
# Row 48 - FirstDigitalUSD
$ws.Range("E48").Value = "  +0.32%  "
# Row 49 - THORChain
Set-TextValue $ws.Range("D49") "8.65"
$ws.Range("E49").Value = "  -5.40%  "

# Row 50 - FLOKI
$ws.Range("E50").Value = "  +0.93%  "
# Row 51 - Monero
Set-TextValue $ws.Range("D51") "129.86"
$ws.Range("E51").Value = "  -4.39%  "

